$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(142).Insert()

$ws.Cells.Item(142, 1).Value = 10
$ws.Cells.Item(142, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(142, 3).Value = "La Araucanía"
$ws.Cells.Item(142, 4).Value = 44463
$ws.Cells.Item(142, 5).Value = 9
$ws.Cells.Item(142, 6).Value = 100112009
$ws.Cells.Item(142, 7).Value = "Acelga"
$ws.Cells.Item(142, 8).Value = "Sin especificar"
$ws.Cells.Item(142, 9).Value = "Primera"
$ws.Cells.Item(142, 10).Value = 30
$ws.Cells.Item(142, 11).Value = 8000
$ws.Cells.Item(142, 12).Value = 8000
$ws.Cells.Item(142, 13).Value = 8000
$ws.Cells.Item(142, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(142, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(142, 16).Value = 667
$ws.Cells.Item(142, 17).Value = 12
$ws.Cells.Item(142, 18).Value = "Hortaliza"

Write-Output "done"
